$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a horizontal run of values starting at column D (4) for a
# given row, one cell at a time (Range.Value array-assignment isn't
# supported by this host, so Cells.Item(...) in a loop is used instead).
function Set-RowValues($row, $values) {
    $col = 4   # column D
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# ---------------------------------------------------------------------------
# Row 2 (2014/12  IFRS연결) : D2:AJ2
# ---------------------------------------------------------------------------
Set-RowValues 2 @(
    9730, 67, 67, 92, 24, 1, 24, 6567, 3938, 2629, 2079, 550, 304, 127, -581, 421,
    561, -434, 1754, 0.68, 0.25, 0.03, 0.39, 149.83, 569.53, 2, 1629.03, 6826, 0.51,
    50, 0.29, 2354.09, 30450420
)

# ---------------------------------------------------------------------------
# Row 3 (2015/12  IFRS연결) : D3:AJ3
# ---------------------------------------------------------------------------
Set-RowValues 3 @(
    11325, 118, 118, 175, 143, 121, 22, 7171, 4416, 2754, 2182, 572, 304, 468, -501, -31,
    517, -48, 1805, 1.04, 1.26, 5.66, 2.08, 160.35, 599.9400000000001, 396, 7.24, 7165, 0.4,
    0, 0, 0, 30450420
)

# ---------------------------------------------------------------------------
# Row 4 (2016/12  IFRS연결) : D4:AJ4
# ---------------------------------------------------------------------------
Set-RowValues 4 @(
    11365, 129, 114, 163, 137, 129, 8, 6897, 4400, 2496, 2317, 179, 304, 413, -352, 211,
    408, 5, 2089, 1.13, 1.21, 5.75, 1.95, 176.28, 643.66, 425, 6.4, 7610, 0.36,
    50, 0.37, 11.76, 30450420
)

# ---------------------------------------------------------------------------
# Row 5 (2017/12  IFRS연결) : D5:AJ5
# ---------------------------------------------------------------------------
Set-RowValues 5 @(
    7210, -171, -171, -206, -236, -216, -20, 6325, 4098, 2228, 2080, 147, 304, 261, -253, -58,
    322, -61, 1998, -2.37, -3.27, -9.81, -3.57, 183.94, 574.37, -708, -2.8, 6832, 0.29,
    0, 0, 0, 30450420
)

# ---------------------------------------------------------------------------
# Row 6 (2018/12  IFRS연결) : D6 onward. Columns J6 and O6 are already blank
# in the source (CAPEX / 자본총계(비지배) were never populated for this row)
# and stay blank. AG6 / AH6 (현금DPS / 현금배당수익률) are cleared entirely
# as part of this edit, so they are handled individually below rather than
# through the generic column-by-column helper.
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = 7937
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = -10
$ws.Range("I6").Value = 15
$ws.Range("K6").Value = 6406
$ws.Range("L6").Value = 4202
$ws.Range("M6").Value = 2204
$ws.Range("N6").Value = 2081
$ws.Range("P6").Value = 304
$ws.Range("Q6").Value = 420
$ws.Range("R6").Value = -131
$ws.Range("S6").Value = -433
$ws.Range("T6").Value = 296
$ws.Range("U6").Value = 124
$ws.Range("V6").Value = 1662
$ws.Range("W6").Value = 0.25
$ws.Range("X6").Value = -0.13
$ws.Range("Y6").Value = 0.72
$ws.Range("Z6").Value = -0.16
$ws.Range("AA6").Value = 190.69
$ws.Range("AB6").Value = 572.97
$ws.Range("AC6").Value = 49
$ws.Range("AD6").Value = 27.58
$ws.Range("AE6").Value = 6833
$ws.Range("AF6").Value = 0.2
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 30450420

# ---------------------------------------------------------------------------
# Rows 7, 8, 9 (2019/12(E), 2020/12(E), 2021/12(E)) lose all of their
# estimate data (columns D through AI); only the row index / label columns
# A-C remain populated.
# ---------------------------------------------------------------------------
$ws.Range("D7:AI9").ClearContents()
